$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF (58) holds a game-date string per team row. It was entered as
# "2-17-2012-13" (day-month / season mash-up) which is ambiguous; correct
# it to the unambiguous ISO form "2013-02-17" for every data row (2-31).
# Force the cell format to Text first so Excel keeps the corrected value
# as a literal string instead of re-parsing it as a date serial number.
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)
    if ($cell.Value2 -eq "2-17-2012-13") {
        $cell.NumberFormat = "@"
        $cell.Value = "2013-02-17"
    }
}
